$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster order/content: "Trayce Jackson-Davis" is replaced by "Keon Johnson",
# and the player rows are reordered.
$players = @(
    @{ Name = "Dejounte Murray";    Pos = "PG,SG";   Team = "New Orleans Pelicans" },
    @{ Name = "Keon Johnson";       Pos = "PG,SG";   Team = "Brooklyn Nets" },
    @{ Name = "Chris Paul";         Pos = "PG";      Team = "San Antonio Spurs" },
    @{ Name = "Tyus Jones";         Pos = "PG";      Team = "Phoenix Suns" },
    @{ Name = "Jaylen Brown";       Pos = "SG,SF";   Team = "Boston Celtics" },
    @{ Name = "Pascal Siakam";      Pos = "SF,PF,C"; Team = "Indiana Pacers" },
    @{ Name = "Khris Middleton";    Pos = "SF";      Team = "Milwaukee Bucks" },
    @{ Name = "Deni Avdija";        Pos = "SF,PF";   Team = "Portland Trail Blazers" },
    @{ Name = "Rudy Gobert";        Pos = "C";       Team = "Minnesota Timberwolves" },
    @{ Name = "Jakob Poeltl";       Pos = "C";       Team = "Toronto Raptors" },
    @{ Name = "Jonas Valanciunas";  Pos = "C";       Team = "Washington Wizards" },
    @{ Name = "Jalen Green";        Pos = "PG,SG";   Team = "Houston Rockets" },
    @{ Name = "Draymond Green";     Pos = "PF,C";    Team = "Golden State Warriors" },
    @{ Name = "Nikola Jokic";       Pos = "C";       Team = "Denver Nuggets" },
    @{ Name = "Russell Westbrook";  Pos = "PG,SG";   Team = "Denver Nuggets" },
    @{ Name = "Paolo Banchero";     Pos = "SF,PF";   Team = "Orlando Magic" },
    @{ Name = "Chet Holmgren";      Pos = "PF,C";    Team = "Oklahoma City Thunder" },
    @{ Name = "Jalen Suggs";        Pos = "PG,SG";   Team = "Orlando Magic" }
)

$row = 2
foreach ($p in $players) {
    $ws.Cells.Item($row, 1).Value = $p.Name
    $ws.Cells.Item($row, 2).Value = $p.Pos
    $ws.Cells.Item($row, 3).Value = $p.Team
    $row++
}
